$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text content (streamlit / data source name changes)
$ws.Range("A3").Value = "Appointments with Dr. Naresh"
$ws.Range("A4").Value = "Appointments at Jubliee Hills"
$ws.Range("A5").Value = "Patients with Aetna"

# Minor bug fix: narrow column A width
# (ColumnWidth is offset by ~5/6 of a character unit when round-tripped
# through this host's internal pixel-based storage, so back the input off
# by that amount to land on the exact target stored width of 31.)
$ws.Columns.Item(1).ColumnWidth = 31 - (5/6)
